$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 22: new activity entry dated 2012-09-14 (Excel serial 41166)
$ws.Range("A22").Value = Get-Date -Year 2012 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B22").Value = "Fixed an issue in the OpenCL device information query methods"

# Update the active selection to B25, matching the author's final cursor position
$ws.Range("B25").Select()
